# Weekly data refresh: insert a new "Choclo" (Dulce o Americano, Primera)
# price record for Región de Arica y Parinacota at row 824, pushing the
# existing rows 824:905 down to 825:906.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 824 (shifts 824:905 -> 825:906).
$ws.Rows(824).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A824").Value = 8
$ws.Range("B824").Value = "Terminal La Palmera de La Serena"
$ws.Range("C824").Value = "Coquimbo"
$ws.Range("D824").Value = 45166
$ws.Range("E824").Value = 4
$ws.Range("F824").Value = 100112024
$ws.Range("G824").Value = "Choclo"
$ws.Range("H824").Value = "Dulce o Americano"
$ws.Range("I824").Value = "Primera"
$ws.Range("J824").Value = 500
$ws.Range("K824").Value = 43000
$ws.Range("L824").Value = 44000
$ws.Range("M824").Value = 43500
$ws.Range("N824").Value = "$/malla 70 unidades"
$ws.Range("O824").Value = "Región de Arica y Parinacota"
$ws.Range("P824").Value = 621
$ws.Range("Q824").Value = 70
$ws.Range("R824").Value = "Hortaliza"
